$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace response items "4)".."9)" under Ref A, comment 6/7 area:
#    - "IV.C" -> "IV.E" in items 4 and 5 (matching the run-split left by Word)
#    - insert ", compared to SQL" in item 6
#    - add the brand-new item 7 paragraph about photon loss errors
#    - item 8 becomes a minimal ")" paragraph that now carries the _GoBack
#      bookmark
#    - item 9 keeps the old content that used to sit right after two empty
#      "7.)"/"8.)" placeholder paragraphs (now removed)
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t.StartsWith("4) We added a paragraph about the considerations")) {
        $startPara = $i
    }
    if ($startPara -ne $null -and $t -like "*connected by quantum entanglement.*") {
        $endPara = $i
        break
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    Write-Host "ERROR: could not locate target paragraph range (start=$startPara end=$endPara)"
} else {
    $rngStart = $d.Paragraphs.Item($startPara).Range.Start
    $rngEnd = $d.Paragraphs.Item($endPara).Range.End
    $rng = $d.Range($rngStart, $rngEnd)
    $xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>4</w:t></w:r><w:r><w:t>) We added a paragraph about the considerations of the phase matching condition in the presence of an optical c</w:t></w:r><w:r><w:t>avity to the end of section IV.E</w:t></w:r><w:r><w:t xml:space="preserve"> in SI.</w:t></w:r></w:p><w:p><w:r><w:t>5</w:t></w:r><w:r><w:t>) We added a fig</w:t></w:r><w:r><w:t>ure (now Fig 1.) to section IV.E</w:t></w:r><w:r><w:t xml:space="preserve"> of SI, which illustrates the orientation of the coherent driving fields with respect to the optical cavity field.</w:t></w:r></w:p><w:p><w:r><w:t>6</w:t></w:r><w:r><w:t xml:space="preserve">) Comparing our results to the standard quantum </w:t></w:r><w:r><w:t>limit</w:t></w:r><w:r><w:t>, as a benchmark, has the advantage of being easily comparable to other results that also compare themselves to SQL. In case of the results reported in Ref [10], they report a 70-fold increase in accuracy of phase measurement, which ideally would translate to the same enhancement in clock stability. We found a 12-fold enhancement in our analysis</w:t></w:r><w:r><w:t>, compared to SQL</w:t></w:r><w:r><w:t>. To make the comparison easier, we changed the sentence about Ref [10] in the introduction to “</w:t></w:r><w:r><w:t>Significant noise reduction has recently been</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>demonstrated with spin-squeezed states in a single ensemble of atoms in</w:t></w:r><w:r><w:t xml:space="preserve"> [10]</w:t></w:r><w:r><w:t>, which reported a 70-fold enhancement of phase measurement</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>accuracy beyond the standard quantum limit.</w:t></w:r><w:r><w:t>”</w:t></w:r></w:p><w:p><w:r><w:t>7</w:t></w:r><w:r><w:t>) We added three paragraphs describing the limitations of our scheme arising from photon loss errors.</w:t></w:r><w:r><w:t xml:space="preserve"> We derived typical maxi</w:t></w:r><w:r><w:t>mal distances, for which the photon propagation loss is not significantly larger than the inherent probabilistic “loss” of the two-photon scheme. We report results for both optical fiber</w:t></w:r><w:r><w:t xml:space="preserve"> links between terrestrial labs and </w:t></w:r><w:r><w:t>free-space optical links between satellites.</w:t></w:r></w:p><w:p><w:r><w:t>8</w:t></w:r><w:r><w:t>)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>9</w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t xml:space="preserve">Yes, local oscillators of the clocks have to be phase locked prior to entangling the atoms. </w:t></w:r><w:r><w:t xml:space="preserve">We added </w:t></w:r><w:r><w:t xml:space="preserve">a </w:t></w:r><w:r><w:t>clarification</w:t></w:r><w:r><w:t xml:space="preserve"> to the introduction: “…</w:t></w:r><w:r><w:t xml:space="preserve"> network of atomic clocks</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>can result in substantial boost of the overall precision if multiple</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>clocks are phase locked and connected by quantum entanglement.</w:t></w:r><w:r><w:t>”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xmlFrag1)
}

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker: it used to render right
#    before "- We moved the lower indices..." but now (because several
#    paragraphs of text were added earlier in the document) it falls right
#    before "- We modified the following sentence" instead.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*We moved the lower indices*") {
        $pr = $p.Range
        $runRng = $d.Range($pr.Start, $pr.End - 1)
        $xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>- We moved the lower indices inside the kets in Eq. 4, so that the description in following text is easier to follow.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $runRng.InsertXML($xmlFrag2)
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*We modified the following sentence*") {
        $pr = $p.Range
        $runRng = $d.Range($pr.Start, $pr.End - 1)
        $xmlFrag3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>- We modified the following sentence</w:t></w:r><w:r><w:br/><w:t>“, which promotes any</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">population in s to r_2, which then blocks the path </w:t></w:r><w:r><w:t xml:space="preserve">via r_1.” to </w:t></w:r><w:r><w:br/><w:t>“</w:t></w:r><w:r><w:t>T</w:t></w:r><w:r><w:t>his promotes any population in s to r_2</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> which then blocks the path </w:t></w:r><w:r><w:t>g</w:t></w:r><w:r><w:t xml:space="preserve"> ↔</w:t></w:r><w:r><w:t xml:space="preserve"> r_1 ↔</w:t></w:r><w:r><w:t xml:space="preserve"> f</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $runRng.InsertXML($xmlFrag3)
        break
    }
}

Write-Host "Done."
